# Add a "Save" column (column H) to the s_vals worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header cell (G1) onto the new
# header cell (H1) so it picks up the same bold/bordered/centered style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Header label for the new column.
$ws.Range("H1").Value = "Save"

# Data values for the new column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0

Write-Host "Added Save column (H1:H3)"
